# Update Ntng2-Lrrc4 LR-pair stats per Dr Hou's advice:
# ligand-/receptor-expressing cell counts (E, K) go from 1 to 3,
# and all dependent expression/specificity metrics are refreshed to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 14.96931633333333
$ws.Range("H2").Value = 44.907949
$ws.Range("I2").Value = 0.521301883166304
$ws.Range("J2").Value = 0.521301883166304
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.956596
$ws.Range("N2").Value = 5.869788
$ws.Range("O2").Value = 0.2233841747945733
$ws.Range("P2").Value = 0.2233841747945733
$ws.Range("Q2").Value = 29.28890446053467
$ws.Range("R2").Value = 263.600140144812
$ws.Range("S2").Value = 0.1164505909899619
$ws.Range("T2").Value = 0.1164505909899619

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 14.96931633333333
$ws.Range("H3").Value = 44.907949
$ws.Range("I3").Value = 0.521301883166304
$ws.Range("J3").Value = 0.521301883166304
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.312913666666667
$ws.Range("N3").Value = 9.938741
$ws.Range("O3").Value = 0.3782346920846191
$ws.Range("P3").Value = 0.3782346920846191
$ws.Range("Q3").Value = 49.59205266135655
$ws.Range("R3").Value = 446.328473952209
$ws.Range("S3").Value = 0.1971744572625391
$ws.Range("T3").Value = 0.1971744572625391

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 14.96931633333333
$ws.Range("H4").Value = 44.907949
$ws.Range("I4").Value = 0.521301883166304
$ws.Range("J4").Value = 0.521301883166304
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.223766666666667
$ws.Range("N4").Value = 6.6713
$ws.Range("O4").Value = 0.253886996482162
$ws.Range("P4").Value = 0.2538869964821621
$ws.Range("Q4").Value = 33.28826668485556
$ws.Range("R4").Value = 299.5944001637
$ws.Range("S4").Value = 0.1323517693775879
$ws.Range("T4").Value = 0.1323517693775879

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 14.96931633333333
$ws.Range("H5").Value = 44.907949
$ws.Range("I5").Value = 0.521301883166304
$ws.Range("J5").Value = 0.521301883166304
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.265607333333333
$ws.Range("N5").Value = 3.796822
$ws.Range("O5").Value = 0.1444941366386454
$ws.Range("P5").Value = 0.1444941366386455
$ws.Range("Q5").Value = 18.94527652645311
$ws.Range("R5").Value = 170.507488738078
$ws.Range("S5").Value = 0.07532506553621511
$ws.Range("T5").Value = 0.07532506553621512

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.459653666666667
$ws.Range("H6").Value = 22.378961
$ws.Range("I6").Value = 0.2597801674844976
$ws.Range("J6").Value = 0.2597801674844975
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.956596
$ws.Range("N6").Value = 5.869788
$ws.Range("O6").Value = 0.2233841747945733
$ws.Range("P6").Value = 0.2233841747945733
$ws.Range("Q6").Value = 14.59552852558533
$ws.Range("R6").Value = 131.359756730268
$ws.Range("S6").Value = 0.05803077834152053
$ws.Range("T6").Value = 0.05803077834152052

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.459653666666667
$ws.Range("H7").Value = 22.378961
$ws.Range("I7").Value = 0.2597801674844976
$ws.Range("J7").Value = 0.2597801674844975
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.312913666666667
$ws.Range("N7").Value = 9.938741
$ws.Range("O7").Value = 0.3782346920846191
$ws.Range("P7").Value = 0.3782346920846191
$ws.Range("Q7").Value = 24.71318858090011
$ws.Range("R7").Value = 222.418697228101
$ws.Range("S7").Value = 0.09825787165818972
$ws.Range("T7").Value = 0.09825787165818969

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.459653666666667
$ws.Range("H8").Value = 22.378961
$ws.Range("I8").Value = 0.2597801674844976
$ws.Range("J8").Value = 0.2597801674844975
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.223766666666667
$ws.Range("N8").Value = 6.6713
$ws.Range("O8").Value = 0.253886996482162
$ws.Range("P8").Value = 0.2538869964821621
$ws.Range("Q8").Value = 16.58852916881111
$ws.Range("R8").Value = 149.2967625193
$ws.Range("S8").Value = 0.06595480646827209
$ws.Range("T8").Value = 0.06595480646827209

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.459653666666667
$ws.Range("H9").Value = 22.378961
$ws.Range("I9").Value = 0.2597801674844976
$ws.Range("J9").Value = 0.2597801674844975
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.265607333333333
$ws.Range("N9").Value = 3.796822
$ws.Range("O9").Value = 0.1444941366386454
$ws.Range("P9").Value = 0.1444941366386455
$ws.Range("Q9").Value = 9.440992384660223
$ws.Range("R9").Value = 84.96893146194201
$ws.Range("S9").Value = 0.03753671101651519
$ws.Range("T9").Value = 0.03753671101651519

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.657723666666667
$ws.Range("H10").Value = 7.973171000000001
$ws.Range("I10").Value = 0.09255441741743681
$ws.Range("J10").Value = 0.09255441741743679
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.956596
$ws.Range("N10").Value = 5.869788
$ws.Range("O10").Value = 0.2233841747945733
$ws.Range("P10").Value = 0.2233841747945733
$ws.Range("Q10").Value = 5.200091495305334
$ws.Range("R10").Value = 46.800823457748
$ws.Range("S10").Value = 0.0206751921583866
$ws.Range("T10").Value = 0.0206751921583866

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.657723666666667
$ws.Range("H11").Value = 7.973171000000001
$ws.Range("I11").Value = 0.09255441741743681
$ws.Range("J11").Value = 0.09255441741743679
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.312913666666667
$ws.Range("N11").Value = 9.938741
$ws.Range("O11").Value = 0.3782346920846191
$ws.Range("P11").Value = 0.3782346920846191
$ws.Range("Q11").Value = 8.804809057523446
$ws.Range("R11").Value = 79.24328151771101
$ws.Range("S11").Value = 0.03500729157295552
$ws.Range("T11").Value = 0.03500729157295551

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.657723666666667
$ws.Range("H12").Value = 7.973171000000001
$ws.Range("I12").Value = 0.09255441741743681
$ws.Range("J12").Value = 0.09255441741743679
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 2.223766666666667
$ws.Range("N12").Value = 6.6713
$ws.Range("O12").Value = 0.253886996482162
$ws.Range("P12").Value = 0.2538869964821621
$ws.Range("Q12").Value = 5.910157299144445
$ws.Range("R12").Value = 53.1914156923
$ws.Range("S12").Value = 0.02349836304926933
$ws.Range("T12").Value = 0.02349836304926933

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.657723666666667
$ws.Range("H13").Value = 7.973171000000001
$ws.Range("I13").Value = 0.09255441741743681
$ws.Range("J13").Value = 0.09255441741743679
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.265607333333333
$ws.Range("N13").Value = 3.796822
$ws.Range("O13").Value = 0.1444941366386454
$ws.Range("P13").Value = 0.1444941366386455
$ws.Range("Q13").Value = 3.363634562506889
$ws.Range("R13").Value = 30.272711062562
$ws.Range("S13").Value = 0.01337357063682534
$ws.Range("T13").Value = 0.01337357063682534

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.628561
$ws.Range("H14").Value = 10.885683
$ws.Range("I14").Value = 0.1263635319317616
$ws.Range("J14").Value = 0.1263635319317616
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.956596
$ws.Range("N14").Value = 5.869788
$ws.Range("O14").Value = 0.2233841747945733
$ws.Range("P14").Value = 0.2233841747945733
$ws.Range("Q14").Value = 7.099627938356
$ws.Range("R14").Value = 63.896651445204
$ws.Range("S14").Value = 0.02822761330470428
$ws.Range("T14").Value = 0.02822761330470428

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.628561
$ws.Range("H15").Value = 10.885683
$ws.Range("I15").Value = 0.1263635319317616
$ws.Range("J15").Value = 0.1263635319317616
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.312913666666667
$ws.Range("N15").Value = 9.938741
$ws.Range("O15").Value = 0.3782346920846191
$ws.Range("P15").Value = 0.3782346920846191
$ws.Range("Q15").Value = 12.02110932723367
$ws.Range("R15").Value = 108.189983945103
$ws.Range("S15").Value = 0.04779507159093479
$ws.Range("T15").Value = 0.04779507159093478

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.628561
$ws.Range("H16").Value = 10.885683
$ws.Range("I16").Value = 0.1263635319317616
$ws.Range("J16").Value = 0.1263635319317616
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 2.223766666666667
$ws.Range("N16").Value = 6.6713
$ws.Range("O16").Value = 0.253886996482162
$ws.Range("P16").Value = 0.2538869964821621
$ws.Range("Q16").Value = 8.069072999766666
$ws.Range("R16").Value = 72.6216569979
$ws.Range("S16").Value = 0.03208205758703273
$ws.Range("T16").Value = 0.03208205758703273

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.628561
$ws.Range("H17").Value = 10.885683
$ws.Range("I17").Value = 0.1263635319317616
$ws.Range("J17").Value = 0.1263635319317616
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.223766666666667
$ws.Range("N17").Value = 6.6713
$ws.Range("O17").Value = 0.1444941366386454
$ws.Range("P17").Value = 0.1444941366386455
$ws.Range("Q17").Value = 4.592333411047333
$ws.Range("R17").Value = 41.331000699426
$ws.Range("S17").Value = 0.0182587894490898
$ws.Range("T17").Value = 0.0182587894490898
